$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 8.644702000000001
    "H2" = 25.934106
    "I2" = 0.2659466972557785
    "J2" = 0.2659466972557785
    "M2" = 164.855011
    "N2" = 494.565033
    "O2" = 0.8897130830256036
    "P2" = 0.8897130830256037
    "Q2" = 1425.122443301722
    "R2" = 12826.1019897155
    "S2" = 0.2366162559359155
    "T2" = 0.2366162559359155
    "G3" = 8.644702000000001
    "H3" = 25.934106
    "I3" = 0.2659466972557785
    "J3" = 0.2659466972557785
    "O3" = 0.009319291703988999
    "P3" = 0.009319291703989001
    "Q3" = 14.92743224351133
    "R3" = 134.346890191602
    "S3" = 0.00247843484943905
    "T3" = 0.002478434849439051
    "G4" = 8.644702000000001
    "H4" = 25.934106
    "I4" = 0.2659466972557785
    "J4" = 0.2659466972557785
    "M4" = 4.66298
    "N4" = 13.98894
    "O4" = 0.02516583685701085
    "P4" = 0.02516583685701086
    "Q4" = 40.31007253196
    "R4" = 362.79065278764
    "S4" = 0.006692771195799777
    "T4" = 0.006692771195799778
    "G5" = 8.644702000000001
    "H5" = 25.934106
    "I5" = 0.2659466972557785
    "J5" = 0.2659466972557785
    "M5" = 14.04531966666667
    "N5" = 42.135959
    "O5" = 0.07580178841339644
    "P5" = 0.07580178841339645
    "Q5" = 121.4176030130727
    "R5" = 1092.758427117654
    "S5" = 0.02015923527462412
    "T5" = 0.02015923527462412
    "I6" = 0.331579210752513
    "J6" = 0.3315792107525131
    "M6" = 164.855011
    "N6" = 494.565033
    "O6" = 0.8897130830256036
    "P6" = 0.8897130830256037
    "Q6" = 1776.825882222573
    "R6" = 15991.43294000315
    "S6" = 0.2950103618658147
    "T6" = 0.2950103618658148
    "I7" = 0.331579210752513
    "J7" = 0.3315792107525131
    "O7" = 0.009319291703988999
    "P7" = 0.009319291703989001
    "S7" = 0.003090083387981115
    "T7" = 0.003090083387981116
    "I8" = 0.331579210752513
    "J8" = 0.3315792107525131
    "M8" = 4.66298
    "N8" = 13.98894
    "O8" = 0.02516583685701085
    "P8" = 0.02516583685701086
    "Q8" = 50.25812380241334
    "R8" = 452.32311422172
    "S8" = 0.00834446832297416
    "T8" = 0.008344468322974164
    "I9" = 0.331579210752513
    "J9" = 0.3315792107525131
    "M9" = 14.04531966666667
    "N9" = 42.135959
    "O9" = 0.07580178841339644
    "P9" = 0.07580178841339645
    "Q9" = 151.3820378066824
    "R9" = 1362.438340260142
    "S9" = 0.02513429717574298
    "T9" = 0.02513429717574299
    "G10" = 7.273908666666667
    "H10" = 21.821726
    "I10" = 0.2237754391117454
    "J10" = 0.2237754391117454
    "M10" = 164.855011
    "N10" = 494.565033
    "O10" = 0.8897130830256036
    "P10" = 0.8897130830256037
    "Q10" = 1199.140293256329
    "R10" = 10792.26263930696
    "S10" = 0.1990959358375192
    "T10" = 0.1990959358375192
    "G11" = 7.273908666666667
    "H11" = 21.821726
    "I11" = 0.2237754391117454
    "J11" = 0.2237754391117454
    "O11" = 0.009319291703988999
    "P11" = 0.009319291703989001
    "Q11" = 12.56038424079355
    "R11" = 113.043458167142
    "S11" = 0.002085428593270584
    "T11" = 0.002085428593270584
    "G12" = 7.273908666666667
    "H12" = 21.821726
    "I12" = 0.2237754391117454
    "J12" = 0.2237754391117454
    "M12" = 4.66298
    "N12" = 13.98894
    "O12" = 0.02516583685701085
    "P12" = 0.02516583685701086
    "Q12" = 33.91809063449334
    "R12" = 305.26281571044
    "S12" = 0.005631496193292149
    "T12" = 0.005631496193292149
    "G13" = 7.273908666666667
    "H13" = 21.821726
    "I13" = 0.2237754391117454
    "J13" = 0.2237754391117454
    "M13" = 14.04531966666667
    "N13" = 42.135959
    "O13" = 0.07580178841339644
    "P13" = 0.07580178841339645
    "Q13" = 102.1643724494704
    "R13" = 919.4793520452341
    "S13" = 0.0169625784876634
    "T13" = 0.0169625784876634
    "G14" = 5.80867
    "H14" = 17.42601
    "I14" = 0.1786986528799631
    "J14" = 0.1786986528799631
    "M14" = 164.855011
    "N14" = 494.565033
    "O14" = 0.8897130830256036
    "P14" = 0.8897130830256037
    "Q14" = 957.58835674537
    "R14" = 8618.295210708331
    "S14" = 0.1589905293863541
    "T14" = 0.1589905293863542
    "G15" = 5.80867
    "H15" = 17.42601
    "I15" = 0.1786986528799631
    "J15" = 0.1786986528799631
    "O15" = 0.009319291703988999
    "P15" = 0.009319291703989001
    "Q15" = 10.03025064946333
    "R15" = 90.27225584517001
    "S15" = 0.00166534487329825
    "T15" = 0.001665344873298251
    "G16" = 5.80867
    "H16" = 17.42601
    "I16" = 0.1786986528799631
    "J16" = 0.1786986528799631
    "M16" = 4.66298
    "N16" = 13.98894
    "O16" = 0.02516583685701085
    "P16" = 0.02516583685701086
    "Q16" = 27.0857120366
    "R16" = 243.7714083294
    "S16" = 0.004497101144944764
    "T16" = 0.004497101144944764
    "G17" = 5.80867
    "H17" = 17.42601
    "I17" = 0.1786986528799631
    "J17" = 0.1786986528799631
    "M17" = 14.04531966666667
    "N17" = 42.135959
    "O17" = 0.07580178841339644
    "P17" = 0.07580178841339645
    "Q17" = 81.58462698817667
    "R17" = 734.2616428935901
    "S17" = 0.01354567747536594
    "T17" = 0.01354567747536594
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

Write-Output "Updated $($updates.Count) cells"